# Update the Rtn4-Tnfrsf19 LR-pairs sheet with the new TPM-based results.
# Rows 2-13 are refreshed in place and rows 14-19 are newly appended, each
# row holding: Sending cluster, Ligand symbol, Receptor symbol, Target
# cluster, then the 16 numeric NATMI metric columns (E..T).

$ws = $excel.ActiveWorkbook.ActiveSheet

$data = @(
    @("ECs", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 80.60050966666667, 241.801529, 0.17420496858261, 0.17420496858261, 1, 0.3333333333333333, 0.037775, 0.113325, 0.05650904967560589, 0.0565090496756059, 3.044684252658333, 27.402158273925, 0.00984415722337207, 0.009844157223372072),
    @("ECs", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 80.60050966666667, 241.801529, 0.17420496858261, 0.17420496858261, 2, 0.6666666666666666, 0.11008, 0.33024, 0.1646728309276161, 0.1646728309276161, 8.872504104106666, 79.85253693696001, 0.0286868253381548, 0.0286868253381548),
    @("ECs", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 80.60050966666667, 241.801529, 0.17420496858261, 0.17420496858261, 3, 1, 0.520622, 1.561866, 0.7788181193967781, 0.7788181193967781, 41.96239854367933, 377.661586893114, 0.1356739860210831, 0.1356739860210831),
    @("FAPs", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 93.34790299999999, 280.043709, 0.2017563980255169, 0.2017563980255169, 1, 0.3333333333333333, 0.037775, 0.113325, 0.05650904967560589, 0.0565090496756059, 3.526217035824999, 31.735953322425, 0.01140106231839525, 0.01140106231839525),
    @("FAPs", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 93.34790299999999, 280.043709, 0.2017563980255169, 0.2017563980255169, 2, 0.6666666666666666, 0.11008, 0.33024, 0.1646728309276161, 0.1646728309276161, 10.27573716224, 92.48163446015998, 0.03322379722062075, 0.03322379722062075),
    @("FAPs", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 93.34790299999999, 280.043709, 0.2017563980255169, 0.2017563980255169, 3, 1, 0.520622, 1.561866, 0.7788181193967781, 0.7788181193967781, 48.59897195566599, 437.390747600994, 0.1571315384865009, 0.1571315384865009),
    @("Inflammatory-Mac", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 82.28866066666666, 246.865982, 0.177853633995942, 0.177853633995942, 1, 0.3333333333333333, 0.037775, 0.113325, 0.05650904967560589, 0.0565090496756059, 3.108454156683333, 27.97608741015, 0.01005033983846371, 0.01005033983846372),
    @("Inflammatory-Mac", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 82.28866066666666, 246.865982, 0.177853633995942, 0.177853633995942, 2, 0.6666666666666666, 0.11008, 0.33024, 0.1646728309276161, 0.1646728309276161, 9.058335766186666, 81.52502189567998, 0.02928766140087586, 0.02928766140087586),
    @("Inflammatory-Mac", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 82.28866066666666, 246.865982, 0.177853633995942, 0.177853633995942, 3, 1, 0.520622, 1.561866, 0.7788181193967781, 0.7788181193967781, 42.84128709360133, 385.5715838424119, 0.1385156327566024, 0.1385156327566024),
    @("MuSCs", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 45.42364, 136.27092, 0.09817585288024938, 0.09817585288024938, 1, 0.3333333333333333, 0.037775, 0.113325, 0.05650904967560589, 0.0565090496756059, 1.715878001, 15.442902009, 0.005547824147354987, 0.005547824147354988),
    @("MuSCs", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 45.42364, 136.27092, 0.09817585288024938, 0.09817585288024938, 2, 0.6666666666666666, 0.11008, 0.33024, 0.1646728309276161, 0.1646728309276161, 5.0002342912, 45.00210862079999, 0.01616689562252381, 0.01616689562252381),
    @("MuSCs", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 45.42364, 136.27092, 0.09817585288024938, 0.09817585288024938, 3, 1, 0.520622, 1.561866, 0.7788181193967781, 0.7788181193967781, 23.64854630408, 212.83691673672, 0.07646113311037059, 0.07646113311037059),
    @("Neutrophils", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 63.60851399999999, 190.825542, 0.1374795175462663, 0.1374795175462663, 1, 0.3333333333333333, 0.037775, 0.113325, 0.05650904967560589, 0.0565090496756059, 2.402811616349999, 21.62530454715, 0.007768836886400293, 0.007768836886400294),
    @("Neutrophils", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 63.60851399999999, 190.825542, 0.1374795175462663, 0.1374795175462663, 2, 0.6666666666666666, 0.11008, 0.33024, 0.1646728309276161, 0.1646728309276161, 7.002025221119999, 63.01822699007999, 0.02263914134890653, 0.02263914134890653),
    @("Neutrophils", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 63.60851399999999, 190.825542, 0.1374795175462663, 0.1374795175462663, 3, 1, 0.520622, 1.561866, 0.7788181193967781, 0.7788181193967781, 33.115991775708, 298.043925981372, 0.1070715393109595, 0.1070715393109595),
    @("Resolving-Mac", "Rtn4", "Tnfrsf19", "ECs", 3, 1, 97.40706899999999, 292.221207, 0.2105296289694155, 0.2105296289694155, 1, 0.3333333333333333, 0.037775, 0.113325, 0.05650904967560589, 0.0565090496756059, 3.679552031474999, 33.115968283275, 0.01189682926161958, 0.01189682926161958),
    @("Resolving-Mac", "Rtn4", "Tnfrsf19", "FAPs", 3, 1, 97.40706899999999, 292.221207, 0.2105296289694155, 0.2105296289694155, 2, 0.6666666666666666, 0.11008, 0.33024, 0.1646728309276161, 0.1646728309276161, 10.72257015552, 96.50313139967999, 0.0346685099965343, 0.0346685099965343),
    @("Resolving-Mac", "Rtn4", "Tnfrsf19", "MuSCs", 3, 1, 97.40706899999999, 292.221207, 0.2105296289694155, 0.2105296289694155, 3, 1, 0.520622, 1.561866, 0.7788181193967781, 0.7788181193967781, 50.712263076918, 456.410367692262, 0.1639642897112616, 0.1639642897112616)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
